# Generate Report for Handoff
# Adds two new handoff rows (bae2ec31-... and d186157a-...) to the
# "Overview", "zh-cn" and "de-de" worksheets / tables.

$wb = $excel.ActiveWorkbook

$bae_guid = "bae2ec31-78c6-40ef-80d7-3da6d94a8b8a"
$d18_guid = "d186157a-a3e7-455d-85d7-39b2f75aecbd"

$bae_md       = "$bae_guid.md"
$bae_md_path  = "e2e\$bae_guid.md"
$d18_md       = "$d18_guid.md"
$d18_md_path  = "e2e\$d18_guid.md"

$bae_zhcn_xlf = "$bae_guid.7e0331f4a09e680130fb9e32419b7eee8f345e87.zh-cn.xlf"
$d18_zhcn_xlf = "$d18_guid.02f38d7b0c6d1123e85b9780bea55f646cd22370.zh-cn.xlf"
$bae_dede_xlf = "$bae_guid.7e0331f4a09e680130fb9e32419b7eee8f345e87.de-de.xlf"
$d18_dede_xlf = "$d18_guid.02f38d7b0c6d1123e85b9780bea55f646cd22370.de-de.xlf"

$bae_ho_date = "2016-09-01 10:48:19"
$d18_ho_date = "2016-09-01 10:48:19"
$bae_xlf_date = "2016-09-01 10:48:13"
$d18_xlf_date = "2016-09-01 10:48:13"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

$baeUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/$bae_md"
$d18Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/$d18_md"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")
$tO = $wsO.ListObjects.Item(1)
$tO.ListRows.Add() | Out-Null
$tO.ListRows.Add() | Out-Null

$wsO.Range("A6").Value2 = $bae_md
$wsO.Range("C6").Value2 = ".md"
$wsO.Range("E6").Value2 = "Ready for handoff"
$wsO.Range("F6").Value2 = "Ready for handoff"
$wsO.Range("G6").Value2 = $bae_ho_date
$wsO.Range("G6").NumberFormat = $dateFmt
$wsO.Hyperlinks.Add($wsO.Range("B6"), $baeUrl, "", "", $bae_md_path) | Out-Null

$wsO.Range("A7").Value2 = $d18_md
$wsO.Range("C7").Value2 = ".md"
$wsO.Range("E7").Value2 = "Ready for handoff"
$wsO.Range("F7").Value2 = "Ready for handoff"
$wsO.Range("G7").Value2 = $d18_ho_date
$wsO.Range("G7").NumberFormat = $dateFmt
$wsO.Hyperlinks.Add($wsO.Range("B7"), $d18Url, "", "", $d18_md_path) | Out-Null

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")
$tZ = $wsZ.ListObjects.Item(1)
$tZ.ListRows.Add() | Out-Null
$tZ.ListRows.Add() | Out-Null

$wsZ.Range("B6").Value2 = ".md"
$wsZ.Range("C6").Value2 = "Ready for handoff"
$wsZ.Range("D6").Value2 = "e2e"
$wsZ.Range("E6").Value2 = "ht"
$wsZ.Range("F6").Value2 = "'False"
$wsZ.Range("G6").Value2 = $bae_zhcn_xlf
$wsZ.Range("H6").Value2 = $bae_xlf_date
$wsZ.Range("H6").NumberFormat = $dateFmt
$wsZ.Range("K6").Value2 = "0001-01-01 00:00:00"
$wsZ.Range("K6").NumberFormat = $dateFmt
$wsZ.Range("M6").Value2 = "'True"
$wsZ.Range("O6").Value2 = "'False"
$wsZ.Hyperlinks.Add($wsZ.Range("A6"), $baeUrl, "", "", $bae_md) | Out-Null

$wsZ.Range("B7").Value2 = ".md"
$wsZ.Range("C7").Value2 = "Ready for handoff"
$wsZ.Range("D7").Value2 = "e2e"
$wsZ.Range("E7").Value2 = "ht"
$wsZ.Range("F7").Value2 = "'False"
$wsZ.Range("G7").Value2 = $d18_zhcn_xlf
$wsZ.Range("H7").Value2 = $d18_xlf_date
$wsZ.Range("H7").NumberFormat = $dateFmt
$wsZ.Range("K7").Value2 = "0001-01-01 00:00:00"
$wsZ.Range("K7").NumberFormat = $dateFmt
$wsZ.Range("M7").Value2 = "'True"
$wsZ.Range("O7").Value2 = "'False"
$wsZ.Hyperlinks.Add($wsZ.Range("A7"), $d18Url, "", "", $d18_md) | Out-Null

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")
$tD = $wsD.ListObjects.Item(1)
$tD.ListRows.Add() | Out-Null
$tD.ListRows.Add() | Out-Null

$wsD.Range("B6").Value2 = ".md"
$wsD.Range("C6").Value2 = "Ready for handoff"
$wsD.Range("D6").Value2 = "e2e"
$wsD.Range("E6").Value2 = "ht"
$wsD.Range("F6").Value2 = "'False"
$wsD.Range("G6").Value2 = $bae_dede_xlf
$wsD.Range("H6").Value2 = $bae_ho_date
$wsD.Range("H6").NumberFormat = $dateFmt
$wsD.Range("K6").Value2 = "0001-01-01 00:00:00"
$wsD.Range("K6").NumberFormat = $dateFmt
$wsD.Range("M6").Value2 = "'True"
$wsD.Range("O6").Value2 = "'False"
$wsD.Hyperlinks.Add($wsD.Range("A6"), $baeUrl, "", "", $bae_md) | Out-Null

$wsD.Range("B7").Value2 = ".md"
$wsD.Range("C7").Value2 = "Ready for handoff"
$wsD.Range("D7").Value2 = "e2e"
$wsD.Range("E7").Value2 = "ht"
$wsD.Range("F7").Value2 = "'False"
$wsD.Range("G7").Value2 = $d18_dede_xlf
$wsD.Range("H7").Value2 = $d18_ho_date
$wsD.Range("H7").NumberFormat = $dateFmt
$wsD.Range("K7").Value2 = "0001-01-01 00:00:00"
$wsD.Range("K7").NumberFormat = $dateFmt
$wsD.Range("M7").Value2 = "'True"
$wsD.Range("O7").Value2 = "'False"
$wsD.Hyperlinks.Add($wsD.Range("A7"), $d18Url, "", "", $d18_md) | Out-Null

Write-Output "Handoff rows added to Overview, zh-cn and de-de sheets."
